$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D36").Value = "Various AI model compressions techniques"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/382"

$ws.Range("D45").Value = "RNN Auto-Encoder (RAE)"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/122"

$ws.Range("D50").Value = "알파텐서"
$ws.Range("E50").Value = "http://incredible.egloos.com/7558034"

$ws.Range("D51").Value = "[python] 딕셔너리에서 value가 가장 큰 key 알아내는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EB%94%95%EC%85%94%EB%84%88%EB%A6%AC%EC%97%90%EC%84%9C-value%EA%B0%80-%EA%B0%80%EC%9E%A5-%ED%81%B0-key-%EC%95%8C%EC%95%84%EB%82%B4%EB%8A%94-%EB%B0%A9%EB%B2%95"
